# Insert a new weekly record at row 68 ("Fruta / hortaliza, semanal").
# This shifts every existing data row (old rows 68..199) down by one
# position (to 69..200), which reproduces the observed diff where each
# row's Fecha/Volumen/Precio values equal the row above's previous
# values, and the former last row (199) becomes the new row 200.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 68 downward (old row 68 -> new row 69, ..., old row 199 -> new row 200)
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with this week's record.
$ws.Cells.Item(68, 1).Value = 4
$ws.Cells.Item(68, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(68, 3).Value = "Los Lagos"
$ws.Cells.Item(68, 4).Value = 44536
$ws.Cells.Item(68, 5).Value = 10
$ws.Cells.Item(68, 6).Value = 100112037
$ws.Cells.Item(68, 7).Value = "Cebollín"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 70
$ws.Cells.Item(68, 11).Value = 6000
$ws.Cells.Item(68, 12).Value = 6000
$ws.Cells.Item(68, 13).Value = 6000
$ws.Cells.Item(68, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(68, 15).Value = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value = 167
$ws.Cells.Item(68, 17).Value = 36
$ws.Cells.Item(68, 18).Value = "Hortaliza"
